# edit.ps1 - applies the tracked changes to the "quiz concept" docx:
#  1. Rewrites / splits the first bullet (Analysis/Stochastik -> + Geometrie/Algebra).
#  2. Rewrites bullets 2 and 3 with new feature text.
#  3. Rewrites bullet 4 (access info) and adds two brand-new bullets after it
#     (score +1 bullet, time-limit bullet) that used to be bullets 2 and 4.
#  4. Relocates the two <w:lastRenderedPageBreak/> markers one heading later.
#  5. Moves the _GoBack bookmark from mid-document to the very last (empty)
#     paragraph of the document.

$d = $word.ActiveDocument

function New-PkgXml($innerBodyXml) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $innerBodyXml + '</w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
}

# Locates the (unique) paragraph containing $searchText and overwrites the
# *whole* paragraph (incl. its end-of-paragraph mark) with $innerParaXml.
# Returns the (reseated) Paragraph object.
function Set-ParagraphByText($doc, [string]$searchText, [string]$innerParaXml) {
    $rng = $doc.Content
    $found = $rng.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Paragraph not found for search text: $searchText"
    }
    $para = $rng.Paragraphs(1)
    $para.Range.InsertXML((New-PkgXml $innerParaXml))
    return $para
}

# Overwrites the content of an already-known Paragraph object.
function Set-ParagraphXml($para, [string]$innerParaXml) {
    $para.Range.InsertXML((New-PkgXml $innerParaXml))
}

# ---------------------------------------------------------------------------
# 1) "50 Aufgaben mit je 4 antworten: Analysis und Stochastik aufgaben"
#    -> three runs: "Aufgab" / "en mit je 4 antworten: Analysis, Geometrie/Algebra " / "und Stochastik "
# ---------------------------------------------------------------------------
$p1 = '<w:p w:rsidR="00E13E73" w:rsidRDefault="00E13E73" w:rsidP="00E13E73">' +
      '<w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
      '<w:r><w:t>Aufgab</w:t></w:r>' +
      '<w:r><w:t xml:space="preserve">en mit je 4 antworten: Analysis, Geometrie/Algebra </w:t></w:r>' +
      '<w:r><w:t xml:space="preserve">und Stochastik </w:t></w:r>' +
      '</w:p>'
Set-ParagraphByText $d "50 Aufgaben mit je 4 antworten: Analysis und Stochastik aufgaben" $p1 | Out-Null

# ---------------------------------------------------------------------------
# 2) "Zugänglich für Schüler, wenn es geht komplett zugänglich online für jeden"
#    -> "Auswahl zwischen den Aufgabenbereichen + Mix-Aufgaben aus allen Bereichen"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Zugänglich für Schüler, wenn es geht komplett zugänglich online für jeden",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Auswahl zwischen den Aufgabenbereichen + Mix-Aufgaben aus allen Bereichen", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) "Speichern der User mit Punktzahl, bei richtigen Antwort +1 bei falschen antworten -1"
#    -> "Bei Falscher Antwort Lösungsweg anzeigen + Tipp Button"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Speichern der User mit Punktzahl, bei richtigen Antwort +1 bei falschen antworten -1",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Bei Falscher Antwort Lösungsweg anzeigen + Tipp Button", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4) "Begrenzte Antwortzeit 20sekunden"
#    -> "Zugänglich für Schüler, online über antonheinrich.github.io/quiz"
#    then two brand-new bullets are appended right after it.
# ---------------------------------------------------------------------------
$p4 = '<w:p w:rsidR="00E13E73" w:rsidRDefault="00E13E73" w:rsidP="00E13E73">' +
      '<w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
      '<w:r><w:t xml:space="preserve">Zugänglich für Schüler, </w:t></w:r>' +
      '<w:r><w:t>online über antonheinrich.github.io/</w:t></w:r>' +
      '<w:proofErr w:type="spellStart"/>' +
      '<w:r><w:t>quiz</w:t></w:r>' +
      '<w:proofErr w:type="spellEnd"/>' +
      '</w:p>'
$para4 = Set-ParagraphByText $d "Begrenzte Antwortzeit 20sekunden" $p4

# New bullet: "Speichern der User mit Punktzahl, bei richtigen Antwort +1 "
$para4.Range.InsertParagraphAfter() | Out-Null
$para5 = $d.Paragraphs($para4.Index + 1)
$p5 = '<w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
      '<w:r><w:t xml:space="preserve">Speichern der User mit Punktzahl, bei richtigen Antwort +1 </w:t></w:r>' +
      '</w:p>'
Set-ParagraphXml $para5 $p5

# New bullet: "Begrenzte Antwortzeit 100 sekunden"
$para5.Range.InsertParagraphAfter() | Out-Null
$para6 = $d.Paragraphs($para5.Index + 1)
$p6 = '<w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
      '<w:r><w:t xml:space="preserve">Begrenzte Antwortzeit 100 </w:t></w:r>' +
      '<w:proofErr w:type="spellStart"/>' +
      '<w:r><w:t>sekunden</w:t></w:r>' +
      '<w:proofErr w:type="spellEnd"/>' +
      '</w:p>'
Set-ParagraphXml $para6 $p6

# ---------------------------------------------------------------------------
# 5) Relocate the two <w:lastRenderedPageBreak/> markers:
#    "c) Benutzerverwaltung"                 loses it
#    "Optional: Feedback (...)"              gains it
#    "Schritt 3: Punkte- und Ergebnislogik"  loses it
#    "Startseite, Quizseite, Ergebnisanzeige." gains it
# ---------------------------------------------------------------------------
$pUserverwaltung = '<w:p w:rsidR="004B70F8" w:rsidRPr="004B70F8" w:rsidRDefault="004B70F8" w:rsidP="004B70F8">' +
    '<w:pPr><w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/>' +
    '<w:outlineLvl w:val="3"/>' +
    '<w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:eastAsia="de-DE"/></w:rPr>' +
    '</w:pPr>' +
    '<w:r w:rsidRPr="004B70F8">' +
    '<w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:eastAsia="de-DE"/></w:rPr>' +
    '<w:t>c) Benutzerverwaltung</w:t></w:r></w:p>'
Set-ParagraphByText $d "c) Benutzerverwaltung" $pUserverwaltung | Out-Null

$pOptionalFeedback = '<w:p w:rsidR="004B70F8" w:rsidRPr="004B70F8" w:rsidRDefault="004B70F8" w:rsidP="004B70F8">' +
    '<w:pPr><w:numPr><w:ilvl w:val="1"/><w:numId w:val="3"/></w:numPr>' +
    '<w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/>' +
    '<w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:eastAsia="de-DE"/></w:rPr>' +
    '</w:pPr>' +
    '<w:r w:rsidRPr="004B70F8">' +
    '<w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:eastAsia="de-DE"/></w:rPr>' +
    '<w:lastRenderedPageBreak/>' +
    '<w:t>Optional: Feedback (z. B. „Super gemacht!“ oder „Übe weiter!“).</w:t></w:r></w:p>'
Set-ParagraphByText $d "Optional: Feedback (z. B. „Super gemacht!“ oder „Übe weiter!“)." $pOptionalFeedback | Out-Null

$pSchritt3 = '<w:p w:rsidR="004B70F8" w:rsidRPr="004B70F8" w:rsidRDefault="004B70F8" w:rsidP="004B70F8">' +
    '<w:pPr><w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/>' +
    '<w:outlineLvl w:val="3"/>' +
    '<w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:eastAsia="de-DE"/></w:rPr>' +
    '</w:pPr>' +
    '<w:r w:rsidRPr="004B70F8">' +
    '<w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:eastAsia="de-DE"/></w:rPr>' +
    '<w:t>Schritt 3: Punkte- und Ergebnislogik</w:t></w:r></w:p>'
Set-ParagraphByText $d "Schritt 3: Punkte- und Ergebnislogik" $pSchritt3 | Out-Null

$pStartseite = '<w:p w:rsidR="004B70F8" w:rsidRPr="004B70F8" w:rsidRDefault="004B70F8" w:rsidP="004B70F8">' +
    '<w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="11"/></w:numPr>' +
    '<w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/>' +
    '<w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:eastAsia="de-DE"/></w:rPr>' +
    '</w:pPr>' +
    '<w:r w:rsidRPr="004B70F8">' +
    '<w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:eastAsia="de-DE"/></w:rPr>' +
    '<w:lastRenderedPageBreak/>' +
    '<w:t>Startseite, Quizseite, Ergebnisanzeige.</w:t></w:r></w:p>'
Set-ParagraphByText $d "Startseite, Quizseite, Ergebnisanzeige." $pStartseite | Out-Null

# ---------------------------------------------------------------------------
# 6) Move the _GoBack bookmark from its mid-document spot to the very last
#    (empty) paragraph of the document.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$d.Bookmarks.Add("_GoBack", $lastPara.Range) | Out-Null

Write-Host "Edit complete."
